$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.187.17'
$ws.Range('E2').Value = '  -2.53%  '
$ws.Range('D3').Value = '3.176.75'
$ws.Range('E3').Value = '  -7.56%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''560.64'
$ws.Range('E5').Value = '  -3.82%  '
$ws.Range('D6').Value = '''171.47'
$ws.Range('E6').Value = '  -1.08%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '''0.601'
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('D9').Value = '3.176.40'
$ws.Range('E9').Value = '  -7.45%  '
$ws.Range('E10').Value = '  -5.55%  '
$ws.Range('D11').Value = '''6.61'
$ws.Range('E11').Value = '  -4.43%  '
$ws.Range('E12').Value = '  -3.02%  '
$ws.Range('D13').Value = '3.721.65'
$ws.Range('E13').Value = '  -7.66%  '
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').Value = '''27.48'
$ws.Range('E15').Value = '  -4.05%  '
$ws.Range('D16').Value = '64.167.76'
$ws.Range('E16').Value = '  -2.64%  '
$ws.Range('E17').Value = '  -4.56%  '
$ws.Range('D18').Value = '3.180.66'
$ws.Range('E18').Value = '  -7.53%  '
$ws.Range('E19').Value = '  -4.34%  '
$ws.Range('D20').Value = '''13.07'
$ws.Range('E20').Value = '  -5.33%  '
$ws.Range('D21').Value = '''352.56'
$ws.Range('E21').Value = '  -4.43%  '
$ws.Range('E22').Value = '  -6.20%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = '''69.19'
$ws.Range('E24').Value = '  -4.33%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').Value = '''0.502'
$ws.Range('E25').Value = '  -5.78%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '''0.0000118'
$ws.Range('E26').Value = '  -2.43%  '
$ws.Range('D27').Value = '''9.46'
$ws.Range('E27').Value = '  -2.72%  '
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '''5.61'
$ws.Range('E30').Value = '  -1.93%  '
$ws.Range('D31').Value = '''0.999'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('E32').Value = '  -4.50%  '
$ws.Range('D33').Value = '''22.07'
$ws.Range('E33').Value = '  -6.60%  '
$ws.Range('E34').Value = '  -5.46%  '
$ws.Range('E35').Value = '  -6.31%  '
$ws.Range('D36').Value = '''157.11'
$ws.Range('E36').Value = '  -2.12%  '
$ws.Range('E37').Value = '  -5.76%  '
$ws.Range('D38').Value = '''26.02'
$ws.Range('E38').Value = '  -9.35%  '
$ws.Range('D39').Value = '''0.796'
$ws.Range('E39').Value = '  -9.25%  '
$ws.Range('D40').Value = '''2.53'
$ws.Range('E40').Value = '  -2.65%  '
$ws.Range('E41').Value = '  -4.48%  '
$ws.Range('D42').Value = '2.656.25'
$ws.Range('E42').Value = '  -4.25%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '''4.16'
$ws.Range('E43').Value = '  -6.63%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '''6.03'
$ws.Range('E44').Value = '  -6.73%  '
$ws.Range('E45').Value = '  -4.07%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = '''328.40'
$ws.Range('E46').Value = '  +1.39%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '''38.78'
$ws.Range('E47').Value = '  -3.37%  '
$ws.Range('D48').Value = '''23.66'
$ws.Range('E48').Value = '  -3.09%  '
$ws.Range('E49').Value = '  -6.52%  '
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('E51').Value = '  +0.03%  '
